$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 38, shifting existing rows 38..85 down to 41..88.
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(38).Insert()
$ws.Rows.Item(38).Insert()

# Populate the three newly inserted rows (38, 39, 40) with the new Castle Brite records.
# Columns A-L (Mercado ID .. Calidad) are identical across all three rows.
$rows = @(38, 39, 40)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).Value2 = 3
    $ws.Cells.Item($r, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($r, 3).Value = "Coquimbo"
    $ws.Cells.Item($r, 4).Value2 = 44536
    $ws.Cells.Item($r, 5).Value2 = 5
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value2 = 100103
    $ws.Cells.Item($r, 8).Value = "Frutos de hueso (carozo)"
    $ws.Cells.Item($r, 9).Value2 = 100103003
    $ws.Cells.Item($r, 10).Value = "Damasco"
    $ws.Cells.Item($r, 11).Value = "Castle Brite"
}

# Row 38: Especial
$ws.Cells.Item(38, 12).Value = "Especial"
$ws.Cells.Item(38, 13).Value2 = 68
$ws.Cells.Item(38, 14).Value2 = 17000
$ws.Cells.Item(38, 15).Value2 = 17000
$ws.Cells.Item(38, 16).Value2 = 17000
$ws.Cells.Item(38, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(38, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(38, 19).Value2 = 1133
$ws.Cells.Item(38, 20).Value2 = 15

# Row 39: Primera
$ws.Cells.Item(39, 12).Value = "Primera"
$ws.Cells.Item(39, 13).Value2 = 70
$ws.Cells.Item(39, 14).Value2 = 15000
$ws.Cells.Item(39, 15).Value2 = 15000
$ws.Cells.Item(39, 16).Value2 = 15000
$ws.Cells.Item(39, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(39, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(39, 19).Value2 = 1000
$ws.Cells.Item(39, 20).Value2 = 15

# Row 40: Segunda
$ws.Cells.Item(40, 12).Value = "Segunda"
$ws.Cells.Item(40, 13).Value2 = 50
$ws.Cells.Item(40, 14).Value2 = 12000
$ws.Cells.Item(40, 15).Value2 = 12000
$ws.Cells.Item(40, 16).Value2 = 12000
$ws.Cells.Item(40, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(40, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(40, 19).Value2 = 800
$ws.Cells.Item(40, 20).Value2 = 15
